$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (even_MAG-GUT449.fa) - entire row
$ws.Rows.Item(2).Delete()

# Delete column C (the "max" column) - entire column
$ws.Columns.Item(3).Delete()

# Update values for the new row 2 (was row 3, even_MAG-GUT48805.fa)
$ws.Cells.Item(2, 2).Value = 1501.686874010694

# Update values for the new row 3 (was row 4, even_MAG-GUT49046.fa)
$ws.Cells.Item(3, 2).Value = 1788.278506209192
